$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 10100880
$ws.Range("J88").Value = 1139.8
$ws.Range("L88").Value = 1139.8
$ws.Range("N88").Value = -1951.8
$ws.Range("H91").Value = 10100880
$ws.Range("J91").Value = 1139.8
$ws.Range("L91").Value = 1139.8
$ws.Range("N91").Value = -3947.8
$ws.Range("H96").Value = 355.125
$ws.Range("I96").Value = 305.85715
$ws.Range("J96").Value = 700
$ws.Range("K96").Value = 917.5714499999999
$ws.Range("L96").Value = 2100
$ws.Range("M96").Value = 455.4285500000001
$ws.Range("N96").Value = -4846
$ws.Range("H97").Value = 10199.429
$ws.Range("J97").Value = 10199.429
$ws.Range("L97").Value = 30598.287
$ws.Range("N97").Value = -31590.287
$ws.Range("H100").Value = 1690.8
$ws.Range("I100").Value = 1690.8
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1690.8
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1149.8
$ws.Range("H112").Value = 2698.75
$ws.Range("J112").Value = 2698.75
$ws.Range("L112").Value = 8096.25
$ws.Range("N112").Value = -10312.25
$ws.Range("H115").Value = 1100.2858
$ws.Range("I115").Value = 867
$ws.Range("K115").Value = 2601
$ws.Range("M115").Value = -1034
$ws.Range("H116").Value = 2487.3845
$ws.Range("I116").Value = 2512.2727
$ws.Range("K116").Value = 2512.2727
$ws.Range("M116").Value = 929.7273
$ws.Range("H118").Value = 602.6923
$ws.Range("I118").Value = 602.6923
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 1808.0769
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = -151.0769
$ws.Range("H137").Value = 1664.6562
$ws.Range("I137").Value = 1673.5518
$ws.Range("J137").Value = 1578.6666
$ws.Range("K137").Value = 5020.6554
$ws.Range("L137").Value = 4735.9998
$ws.Range("M137").Value = -2470.6554
$ws.Range("N137").Value = -9835.9998
$ws.Range("H138").Value = 361924.6
$ws.Range("I138").Value = 6421.9165
$ws.Range("J138").Value = 435476.88
$ws.Range("K138").Value = 19265.7495
$ws.Range("L138").Value = 1306430.64
$ws.Range("M138").Value = -14125.7495
$ws.Range("N138").Value = -1316710.64
$ws.Range("N100").ClearContents()
$ws.Range("N118").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3301.6785
$ws.Range("I61").Value = 2092.8096
$ws.Range("J61").Value = 6928.2856
$ws.Range("K61").Value = 2092.8096
$ws.Range("L61").Value = 6928.2856
$ws.Range("M61").Value = -1880.8096
$ws.Range("N61").Value = -7352.2856
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("H97").Value = 1385.7778
$ws.Range("I97").Value = 1414.1875
$ws.Range("J97").Value = 1158.5
$ws.Range("K97").Value = 1414.1875
$ws.Range("L97").Value = 1158.5
$ws.Range("M97").Value = -918.1875
$ws.Range("N97").Value = -2150.5
$ws.Range("H117").Value = 49998.5
$ws.Range("J117").Value = 49998.5
$ws.Range("L117").Value = 49998.5
$ws.Range("N117").Value = -59176.5
$ws.Range("H132").Value = 2223
$ws.Range("I132").Value = 1954.4572
$ws.Range("K132").Value = 5863.3716
$ws.Range("M132").Value = -3333.3716
$ws.Range("H136").Value = 3301.6785
$ws.Range("I136").Value = 2092.8096
$ws.Range("J136").Value = 6928.2856
$ws.Range("K136").Value = 6278.4288
$ws.Range("L136").Value = 20784.8568
$ws.Range("M136").Value = -3728.4288
$ws.Range("N136").Value = -25884.8568
$ws.Range("N96").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 41667276
$ws.Range("J94").Value = 1556.4
$ws.Range("L94").Value = 1556.4
$ws.Range("N94").Value = -2458.4
$ws.Range("H99").Value = 2909.889
$ws.Range("I99").Value = 2095.1365
$ws.Range("K99").Value = 2095.1365
$ws.Range("M99").Value = -597.1365000000001
$ws.Range("H107").Value = 1478.6765
$ws.Range("I107").Value = 1294.0385
$ws.Range("J107").Value = 2078.75
$ws.Range("K107").Value = 1294.0385
$ws.Range("L107").Value = 2078.75
$ws.Range("M107").Value = 625.9614999999999
$ws.Range("N107").Value = -5918.75
$ws.Range("H134").Value = 2734.4722
$ws.Range("I134").Value = 2171.926
$ws.Range("K134").Value = 6515.778
$ws.Range("M134").Value = -3980.778

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2949.6667
$ws.Range("I22").Value = 2299
$ws.Range("K22").Value = 2299
$ws.Range("M22").Value = -1949
$ws.Range("H92").Value = 43117.57
$ws.Range("J92").Value = 43117.57
$ws.Range("L92").Value = 43117.57
$ws.Range("N92").Value = -48109.57
$ws.Range("H93").Value = 4343.5
$ws.Range("I93").Value = 4343.5
$ws.Range("K93").Value = 4343.5
$ws.Range("M93").Value = -2471.5
$ws.Range("H95").Value = 18999.5
$ws.Range("J95").Value = 18999.5
$ws.Range("L95").Value = 18999.5
$ws.Range("N95").Value = -24491.5
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("H120").Value = 37880.832
$ws.Range("J120").Value = 37880.832
$ws.Range("L120").Value = 37880.832
$ws.Range("N120").Value = -45138.832
$ws.Range("H121").Value = 40325.25
$ws.Range("J121").Value = 40325.25
$ws.Range("L121").Value = 40325.25
$ws.Range("N121").Value = -42945.25
$ws.Range("H132").Value = 4360.731
$ws.Range("I132").Value = 3968.25
$ws.Range("J132").Value = 5669
$ws.Range("K132").Value = 11904.75
$ws.Range("L132").Value = 17007
$ws.Range("M132").Value = -9374.75
$ws.Range("N132").Value = -22067
$ws.Range("H141").Value = 390833.94
$ws.Range("J141").Value = 390833.94
$ws.Range("L141").Value = 390833.94
$ws.Range("N141").Value = -401193.94
$ws.Range("N114").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 531.6667
$ws.Range("I92").Value = 531.6667
$ws.Range("K92").Value = 1595.0001
$ws.Range("M92").Value = -347.0001
$ws.Range("H113").Value = 1666.4762
$ws.Range("J113").Value = 2051.5715
$ws.Range("L113").Value = 6154.7145
$ws.Range("N113").Value = -10494.7145
$ws.Range("H114").Value = 1681.6154
$ws.Range("I114").Value = 720.4
$ws.Range("J114").Value = 2282.375
$ws.Range("K114").Value = 2161.2
$ws.Range("L114").Value = 6847.125
$ws.Range("M114").Value = 1092.8
$ws.Range("N114").Value = -13355.125
$ws.Range("H117").Value = 2331.6924
$ws.Range("I117").Value = 2197.5
$ws.Range("K117").Value = 6592.5
$ws.Range("M117").Value = -3150.5
$ws.Range("H123").Value = 2066.7742
$ws.Range("I123").Value = 712.5
$ws.Range("J123").Value = 2391.8
$ws.Range("K123").Value = 2137.5
$ws.Range("L123").Value = 7175.400000000001
$ws.Range("M123").Value = 312.5
$ws.Range("N123").Value = -12075.4
$ws.Range("H139").Value = 2969
$ws.Range("I139").Value = 2582.375
$ws.Range("K139").Value = 7747.125
$ws.Range("M139").Value = -2607.125
$ws.Range("H140").Value = 17272.37
$ws.Range("I140").Value = 12019.556
$ws.Range("K140").Value = 36058.66800000001
$ws.Range("M140").Value = -30878.66800000001
$ws.Range("H141").Value = 35000
$ws.Range("I141").Value = 20000
$ws.Range("K141").Value = 60000
$ws.Range("M141").Value = -54820

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3827
$ws.Range("I132").Value = 3237.6924
$ws.Range("J132").Value = 4728.294
$ws.Range("K132").Value = 9713.0772
$ws.Range("L132").Value = 14184.882
$ws.Range("M132").Value = -7183.0772
$ws.Range("N132").Value = -19244.882
$ws.Range("H133").Value = 201199.2
$ws.Range("J133").Value = 201199.2
$ws.Range("L133").Value = 201199.2
$ws.Range("N133").Value = -211319.2
$ws.Range("H135").Value = 69641.71000000001
$ws.Range("J135").Value = 69641.71000000001
$ws.Range("L135").Value = 69641.71000000001
$ws.Range("N135").Value = -79781.71000000001
$ws.Range("H140").Value = 74999.234
$ws.Range("J140").Value = 74999.234
$ws.Range("L140").Value = 74999.234
$ws.Range("N140").Value = -85359.234

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3575.2646
$ws.Range("J132").Value = 4845.231
$ws.Range("L132").Value = 14535.693
$ws.Range("N132").Value = -19595.693
$ws.Range("H136").Value = 5890.731
$ws.Range("I136").Value = 4610.6
$ws.Range("K136").Value = 13831.8
$ws.Range("M136").Value = -11281.8

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4522.923
$ws.Range("I132").Value = 4246.9707
$ws.Range("K132").Value = 12740.9121
$ws.Range("M132").Value = -10210.9121
